$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update participant/visit counts in column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 7458
$ws1.Range("F15").Value = 2991
$ws1.Range("F18").Value = 692
$ws1.Range("F23").Value = 184
$ws1.Range("F25").Value = 206

# Sheet "全部类型" (All types) - same updates mirrored here
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 7458
$ws4.Range("F19").Value = 2991
$ws4.Range("F23").Value = 692
$ws4.Range("F29").Value = 184
$ws4.Range("F31").Value = 206
